# Apply the updated cryptocurrency market data (price + 1h volume change,
# and a couple of rank swaps that moved coin name/link/price/volume together)
# onto the existing "cryptos" worksheet, cell by cell, per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    # Column D holds numeric-looking strings (e.g. "613.97", "0.999",
    # "0.0000224") that must stay TEXT (matches the source's inlineStr
    # cells) instead of being auto-coerced to Double by the Value setter.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '67.245.61'
$ws.Cells.Item(2, 5).Value = '  -3.10%  '
# Row 3
Set-TextCell 3 4 '3.546.60'
$ws.Cells.Item(3, 5).Value = '  -3.45%  '
# Row 4
Set-TextCell 4 4 '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.14%  '
# Row 5
Set-TextCell 5 4 '613.97'
$ws.Cells.Item(5, 5).Value = '  -4.83%  '
# Row 6
Set-TextCell 6 4 '155.25'
$ws.Cells.Item(6, 5).Value = '  -2.82%  '
# Row 7
Set-TextCell 7 4 '3.542.59'
$ws.Cells.Item(7, 5).Value = '  -3.40%  '
# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.07%  '
# Row 9
Set-TextCell 9 4 '0.486'
$ws.Cells.Item(9, 5).Value = '  -1.84%  '
# Row 10
Set-TextCell 10 4 '0.142'
$ws.Cells.Item(10, 5).Value = '  -2.01%  '
# Row 11
Set-TextCell 11 4 '6.90'
$ws.Cells.Item(11, 5).Value = '  -2.45%  '
# Row 12
Set-TextCell 12 4 '0.432'
$ws.Cells.Item(12, 5).Value = '  -3.88%  '
# Row 13
Set-TextCell 13 4 '0.0000224'
$ws.Cells.Item(13, 5).Value = '  -3.78%  '
# Row 14
Set-TextCell 14 4 '32.21'
$ws.Cells.Item(14, 5).Value = '  -1.75%  '
# Row 15
Set-TextCell 15 4 '4.136.03'
$ws.Cells.Item(15, 5).Value = '  -3.76%  '
# Row 16
Set-TextCell 16 4 '3.529.00'
$ws.Cells.Item(16, 5).Value = '  -3.93%  '
# Row 17
Set-TextCell 17 4 '67.188.30'
$ws.Cells.Item(17, 5).Value = '  -3.21%  '
# Row 19
Set-TextCell 19 4 '6.41'
$ws.Cells.Item(19, 5).Value = '  -1.35%  '
# Row 20
Set-TextCell 20 4 '15.52'
$ws.Cells.Item(20, 5).Value = '  -3.10%  '
# Row 21
Set-TextCell 21 4 '454.40'
$ws.Cells.Item(21, 5).Value = '  -2.58%  '
# Row 22
Set-TextCell 22 4 '9.42'
$ws.Cells.Item(22, 5).Value = '  -4.93%  '
# Row 23
Set-TextCell 23 4 '0.645'
$ws.Cells.Item(23, 5).Value = '  -0.15%  '
# Row 24
Set-TextCell 24 4 '79.09'
$ws.Cells.Item(24, 5).Value = '  -0.48%  '
# Row 25
$ws.Cells.Item(25, 5).Value = '  -0.11%  '
# Row 26
Set-TextCell 26 4 '3.669.44'
$ws.Cells.Item(26, 5).Value = '  -3.96%  '
# Row 27
Set-TextCell 27 4 '0.0000125'
$ws.Cells.Item(27, 5).Value = '  -0.66%  '
# Row 28
Set-TextCell 28 4 '10.47'
$ws.Cells.Item(28, 5).Value = '  -3.82%  '
# Row 29
Set-TextCell 29 4 '8.41'
$ws.Cells.Item(29, 5).Value = '  -7.14%  '
# Row 30
Set-TextCell 30 4 '2.57'
$ws.Cells.Item(30, 5).Value = '  -2.06%  '
# Row 31
Set-TextCell 31 4 '1.70'
$ws.Cells.Item(31, 5).Value = '  -1.34%  '
# Row 32
Set-TextCell 32 4 '0.999'
$ws.Cells.Item(32, 5).Value = '  +0.12%  '
# Row 33
Set-TextCell 33 4 '26.06'
$ws.Cells.Item(33, 5).Value = '  -3.01%  '
# Row 34
Set-TextCell 34 4 '1.91'
$ws.Cells.Item(34, 5).Value = '  -5.04%  '
# Row 35
$ws.Cells.Item(35, 2).Value = 'NEARProtocol'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 35 4 '6.24'
$ws.Cells.Item(35, 5).Value = '  -3.40%  '
# Row 36
$ws.Cells.Item(36, 2).Value = 'Kaspa'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 36 4 '0.159'
$ws.Cells.Item(36, 5).Value = '  -3.75%  '
# Row 37
Set-TextCell 37 4 '3.531.63'
$ws.Cells.Item(37, 5).Value = '  -3.68%  '
# Row 38
Set-TextCell 38 4 '8.11'
$ws.Cells.Item(38, 5).Value = '  -3.75%  '
# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.06%  '
# Row 40
$ws.Cells.Item(40, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell 40 4 '0.996'
$ws.Cells.Item(40, 5).Value = '  -0.47%  '
# Row 41
$ws.Cells.Item(41, 2).Value = 'Monero'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 41 4 '176.27'
$ws.Cells.Item(41, 5).Value = '  -1.18%  '
# Row 42
Set-TextCell 42 4 '5.64'
$ws.Cells.Item(42, 5).Value = '  -4.21%  '
# Row 43
$ws.Cells.Item(43, 2).Value = 'Stacks'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell 43 4 '2.15'
$ws.Cells.Item(43, 5).Value = '  -1.12%  '
# Row 44
$ws.Cells.Item(44, 2).Value = 'Hedera'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 44 4 '0.0882'
$ws.Cells.Item(44, 5).Value = '  -1.98%  '
# Row 45
Set-TextCell 45 4 '0.895'
$ws.Cells.Item(45, 5).Value = '  -3.18%  '
# Row 46
Set-TextCell 46 4 '28.82'
$ws.Cells.Item(46, 5).Value = '  +5.94%  '
# Row 47
Set-TextCell 47 4 '45.77'
$ws.Cells.Item(47, 5).Value = '  -1.86%  '
# Row 48
Set-TextCell 48 4 '2.73'
$ws.Cells.Item(48, 5).Value = '  -0.58%  '
# Row 49
$ws.Cells.Item(49, 2).Value = 'ONDO'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextCell 49 4 '1.23'
$ws.Cells.Item(49, 5).Value = '  -1.60%  '
# Row 50
$ws.Cells.Item(50, 2).Value = 'Cosmos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 50 4 '7.69'
$ws.Cells.Item(50, 5).Value = '  -1.97%  '
# Row 51
$ws.Cells.Item(51, 5).Value = '  -3.26%  '
